$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 173032
$ws.Range("C4").Value = 163785
$ws.Range("C5").Value = 9247
$ws.Range("C6").Value = 530
$ws.Range("C7").Value = 5.34
$ws.Range("C8").Value = 66.18000000000001
